$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark; it will be re-created further
#    down, wrapping "900", in the new content below.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2. Locate the paragraph "Küçük harfe çevirmemem gerek!!!!" (unchanged) and
#    the blank list paragraph that immediately follows it. That blank
#    paragraph is the anchor we will expand into the six new paragraphs
#    described by the commit. The final paragraph of the inserted XML
#    package below (the page-break paragraph) takes over that existing
#    blank paragraph's slot, exactly as happens when Word's COM InsertXML
#    absorbs the host paragraph.
$rng = $d.Content
$rng.Find.Execute("k harfe çevirmemem gerek!!!!", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hostPara = $rng.Paragraphs(1)
$blankPara = $hostPara.Next()
$insertPt = $blankPara.Range
$insertPt.Collapse(1)

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="AralkYok"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Formata uygun olarak</w:t></w:r><w:r><w:t xml:space="preserve"> tagleme işlemi yapıldı.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AralkYok"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>CRF tarafı tamam.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AralkYok"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>CRF API olarak ayağa kaldırıldı.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AralkYok"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Test için 200 adet makale ayrıldı.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AralkYok"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Doğum tarihinde tanımlı olan </w:t></w:r><w:r><w:t>900</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> adet makale işaretlendi.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="AralkYok"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p><w:p><w:r><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPt.InsertXML($xml)

Write-Output "Inserted new progress paragraphs after 'Küçük harfe çevirmemem gerek!!!!'"
